$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing bold/bordered cells onto newly-added cells first
$ws.Range("A2").Copy($ws.Range("A8:A11"))
$ws.Range("V1").Copy($ws.Range("W1"))

# Populate all cell values (header row, labels, and data)
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18
$ws.Range("U1").Value = 19
$ws.Range("V1").Value = 20
$ws.Range("W1").Value = 21
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[1, 1, 1]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[2, 2, 2]"
$ws.Range("H2").Value = "[4, 0, 0]"
$ws.Range("I2").Value = "[3, 3, 1]"
$ws.Range("J2").Value = "[4, 2, 0]"
$ws.Range("K2").Value = "[4, 2, 2]"
$ws.Range("L2").Value = "[5, 1, 1]"
$ws.Range("M2").Value = "[3, 3, 3]"
$ws.Range("N2").Value = "1Pair-A"
$ws.Range("O2").Value = "1Pair-B"
$ws.Range("P2").Value = "2Pairs-A"
$ws.Range("Q2").Value = "2Pairs-B"
$ws.Range("R2").Value = "3Pairs-A"
$ws.Range("S2").Value = "3Pairs-B"
$ws.Range("T2").Value = "3Pairs-C"
$ws.Range("U2").Value = "4Pairs"
$ws.Range("V2").Value = "5A4F"
$ws.Range("W2").Value = "MaxUnique"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Equal Angle"
$ws.Range("C3").Value = 0.9840417867435158
$ws.Range("D3").Value = 0.8732204610951009
$ws.Range("E3").Value = 1.110965417867435
$ws.Range("F3").Value = 0.9556340057636887
$ws.Range("G3").Value = 0.9840417867435158
$ws.Range("H3").Value = 0.8732204610951009
$ws.Range("I3").Value = 1.059466858789625
$ws.Range("J3").Value = 1.004690201729107
$ws.Range("K3").Value = 0.9870317002881844
$ws.Range("L3").Value = 0.9169236311239193
$ws.Range("M3").Value = 0.9840417867435158
$ws.Range("N3").Value = 0.9840417867435158
$ws.Range("O3").Value = 1.110965417867435
$ws.Range("P3").Value = 0.992092939481268
$ws.Range("Q3").Value = 1.033299711815562
$ws.Range("R3").Value = 0.9894092219020173
$ws.Range("S3").Value = 0.9799399615754082
$ws.Range("T3").Value = 0.9894092219020173
$ws.Range("U3").Value = 0.9809654178674352
$ws.Range("V3").Value = 0.9815806916426514
$ws.Range("W3").Value = 0.9864967579250721
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "CLR"
$ws.Range("C4").Value = 1.014622979784832
$ws.Range("D4").Value = 0.9950939054233108
$ws.Range("E4").Value = 0.9870975407402144
$ws.Range("F4").Value = 0.9908578525858583
$ws.Range("G4").Value = 1.014622979784832
$ws.Range("H4").Value = 0.9950939054233108
$ws.Range("I4").Value = 0.9947015877273128
$ws.Range("J4").Value = 0.9884275221326881
$ws.Range("K4").Value = 0.9991745937072521
$ws.Range("L4").Value = 0.9891407001886157
$ws.Range("M4").Value = 1.01460776834583
$ws.Range("N4").Value = 1.014622979784832
$ws.Range("O4").Value = 0.9870975407402144
$ws.Range("P4").Value = 0.9910957230817625
$ws.Range("Q4").Value = 0.9889776966630364
$ws.Range("R4").Value = 0.9989381419827855
$ws.Range("S4").Value = 0.9910164329164611
$ws.Range("T4").Value = 0.9989381419827855
$ws.Range("U4").Value = 0.9969180696335538
$ws.Range("V4").Value = 1.000459051663809
$ws.Range("W4").Value = 0.9948895852862605
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "BT8Hex"
$ws.Range("C5").Value = 1.026991197513424
$ws.Range("D5").Value = 0.9950887042466666
$ws.Range("E5").Value = 0.9833508948227769
$ws.Range("F5").Value = 0.9888510106437166
$ws.Range("G5").Value = 1.026991197513424
$ws.Range("H5").Value = 0.9950887042466666
$ws.Range("I5").Value = 0.9956971765054171
$ws.Range("J5").Value = 0.9834554141958974
$ws.Range("K5").Value = 1.002219417062731
$ws.Range("L5").Value = 0.9864386844527027
$ws.Range("M5").Value = 1.027020112578161
$ws.Range("N5").Value = 1.026991197513424
$ws.Range("O5").Value = 0.9833508948227769
$ws.Range("P5").Value = 0.9892197995347218
$ws.Range("Q5").Value = 0.9861009527332467
$ws.Range("R5").Value = 1.001810265527623
$ws.Range("S5").Value = 0.9890968699043867
$ws.Range("T5").Value = 1.001810265527623
$ws.Range("U5").Value = 0.9985704518066461
$ws.Range("V5").Value = 1.004254600948002
$ws.Range("W5").Value = 0.9952615624304165
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Spiral"
$ws.Range("C6").Value = 0.9941638581492266
$ws.Range("D6").Value = 0.994619251703808
$ws.Range("E6").Value = 0.9962784887745791
$ws.Range("F6").Value = 0.9931141819753077
$ws.Range("G6").Value = 0.9941638581492266
$ws.Range("H6").Value = 0.994619251703808
$ws.Range("I6").Value = 0.9947855173201016
$ws.Range("J6").Value = 0.9961144329414382
$ws.Range("K6").Value = 0.9936797700849543
$ws.Range("L6").Value = 0.9921848445813087
$ws.Range("M6").Value = 0.9941532572641268
$ws.Range("N6").Value = 0.9941638581492266
$ws.Range("O6").Value = 0.9962784887745791
$ws.Range("P6").Value = 0.9954488702391936
$ws.Range("Q6").Value = 0.9946963353749434
$ws.Range("R6").Value = 0.9950205328758712
$ws.Range("S6").Value = 0.9946706408178984
$ws.Range("T6").Value = 0.9950205328758712
$ws.Range("U6").Value = 0.9945439451507303
$ws.Range("V6").Value = 0.9944679277504296
$ws.Range("W6").Value = 0.9943675431913406
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "OffsetF"
$ws.Range("C7").Value = 2.243664783995547
$ws.Range("D7").Value = 0.1435974566754628
$ws.Range("E7").Value = 0.8944693780818429
$ws.Range("F7").Value = 0.7641082459616358
$ws.Range("G7").Value = 2.243664783995547
$ws.Range("H7").Value = 0.1435974566754628
$ws.Range("I7").Value = 1.287896716225035
$ws.Range("J7").Value = 0.5982434962969629
$ws.Range("K7").Value = 1.398067464417403
$ws.Range("L7").Value = 0.3622893947295112
$ws.Range("M7").Value = 2.24313625412751
$ws.Range("N7").Value = 2.243664783995547
$ws.Range("O7").Value = 0.8944693780818429
$ws.Range("P7").Value = 0.5190334173786528
$ws.Range("Q7").Value = 0.8292888120217394
$ws.Range("R7").Value = 1.093910539584284
$ws.Range("S7").Value = 0.6007250269063138
$ws.Range("T7").Value = 1.093910539584284
$ws.Range("U7").Value = 1.011459966178622
$ws.Range("V7").Value = 1.257900929742007
$ws.Range("W7").Value = 0.9615421170479252
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "OffsetA"
$ws.Range("C8").Value = 0.9800356093511665
$ws.Range("D8").Value = 0.5449080945638909
$ws.Range("E8").Value = 1.22692954780446
$ws.Range("F8").Value = 0.9446111129847837
$ws.Range("G8").Value = 0.9800356093511665
$ws.Range("H8").Value = 0.5449080945638909
$ws.Range("I8").Value = 1.154932885282359
$ws.Range("J8").Value = 1.003821204405187
$ws.Range("K8").Value = 1.041304524754317
$ws.Range("L8").Value = 0.7515175661122454
$ws.Range("M8").Value = 0.9800540954002057
$ws.Range("N8").Value = 0.9800356093511665
$ws.Range("O8").Value = 1.22692954780446
$ws.Range("P8").Value = 0.8859188211841756
$ws.Range("Q8").Value = 1.085770330394622
$ws.Range("R8").Value = 0.9172910839065059
$ws.Range("S8").Value = 0.905482918451045
$ws.Range("T8").Value = 0.9172910839065059
$ws.Range("U8").Value = 0.9241210911760754
$ws.Range("V8").Value = 0.9353039948110936
$ws.Range("W8").Value = 0.9560075681573011
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "RD Single"
$ws.Range("C9").Value = 0.07000000000000001
$ws.Range("D9").Value = 3.94
$ws.Range("E9").Value = 0.21
$ws.Range("F9").Value = 1.3
$ws.Range("G9").Value = 0.07000000000000001
$ws.Range("H9").Value = 3.94
$ws.Range("I9").Value = 0.18
$ws.Range("J9").Value = 1.14
$ws.Range("K9").Value = 0.46
$ws.Range("L9").Value = 2.52
$ws.Range("M9").Value = 0.07000000000000001
$ws.Range("N9").Value = 0.07000000000000001
$ws.Range("O9").Value = 0.21
$ws.Range("P9").Value = 2.075
$ws.Range("Q9").Value = 0.755
$ws.Range("R9").Value = 1.406666666666667
$ws.Range("S9").Value = 1.816666666666667
$ws.Range("T9").Value = 1.406666666666667
$ws.Range("U9").Value = 1.38
$ws.Range("V9").Value = 1.118
$ws.Range("W9").Value = 1.2275
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "TD Single"
$ws.Range("C10").Value = 0.67
$ws.Range("D10").Value = 0.21
$ws.Range("E10").Value = 1.97
$ws.Range("F10").Value = 0.7
$ws.Range("G10").Value = 0.67
$ws.Range("H10").Value = 0.21
$ws.Range("I10").Value = 1.52
$ws.Range("J10").Value = 1.16
$ws.Range("K10").Value = 0.86
$ws.Range("L10").Value = 0.46
$ws.Range("M10").Value = 0.67
$ws.Range("N10").Value = 0.67
$ws.Range("O10").Value = 1.97
$ws.Range("P10").Value = 1.09
$ws.Range("Q10").Value = 1.335
$ws.Range("R10").Value = 0.9500000000000001
$ws.Range("S10").Value = 0.96
$ws.Range("T10").Value = 0.9500000000000001
$ws.Range("U10").Value = 0.8875
$ws.Range("V10").Value = 0.844
$ws.Range("W10").Value = 0.94375
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C11").Value = 0.9990330497856222
$ws.Range("D11").Value = 0.9974757325517742
$ws.Range("E11").Value = 0.9915490378394113
$ws.Range("F11").Value = 0.9933549872629314
$ws.Range("G11").Value = 0.9990330497856222
$ws.Range("H11").Value = 0.9974757325517742
$ws.Range("I11").Value = 0.9931845804783355
$ws.Range("J11").Value = 0.9942157476082852
$ws.Range("K11").Value = 0.994935224634625
$ws.Range("L11").Value = 0.9932631714184019
$ws.Range("M11").Value = 0.9990350350090685
$ws.Range("N11").Value = 0.9990330497856222
$ws.Range("O11").Value = 0.9915490378394113
$ws.Range("P11").Value = 0.9945123851955928
$ws.Range("Q11").Value = 0.9924520125511713
$ws.Range("R11").Value = 0.9960192733922693
$ws.Range("S11").Value = 0.9941265858847057
$ws.Range("T11").Value = 0.9960192733922693
$ws.Range("U11").Value = 0.9953532018599348
$ws.Range("V11").Value = 0.9960891714450723
$ws.Range("W11").Value = 0.9946264414474233
